$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add label for median row
$ws.Range("A57").Value = "Median"

# Add median formulas for columns B through E (rows 2-56)
$ws.Range("B57").Formula = "=MEDIAN(B2:B56)"
$ws.Range("C57").Formula = "=MEDIAN(C2:C56)"
$ws.Range("D57").Formula = "=MEDIAN(D2:D56)"
$ws.Range("E57").Formula = "=MEDIAN(E2:E56)"

# Copy number format style from existing data cells (s="2") to the new median values
$ws.Range("B57:E57").NumberFormat = $ws.Range("B2").NumberFormat

# Apply header-like style to A57 (bold, centered, with border) similar to A1-A56 but with left+right border
$ws.Range("A57").Font.Bold = $true
$ws.Range("A57").HorizontalAlignment = -4108
$ws.Range("A57").VerticalAlignment = -4160
$ws.Range("A57").Borders.Item(7).LineStyle = 1
$ws.Range("A57").Borders.Item(7).Weight = 2
$ws.Range("A57").Borders.Item(10).LineStyle = 1
$ws.Range("A57").Borders.Item(10).Weight = 2

$wb.Save()
